$wb = $excel.ActiveWorkbook

# Rename the "Template" sheet to "Authors"
$authors = $wb.Worksheets.Item("Template")
$authors.Name = "Authors"

# Fix typo "Institude" -> "Institute" in the affiliation text (cell F4)
$authors.Range("F4").Value = "Center for Biomedical Informatics and Information Technology, National Cancer Institute, Rockville, MD, USA"

# Update the active selection on the Authors sheet
$authors.Activate() | Out-Null
$authors.Range("A2").Select() | Out-Null
